# Updated cryptos list on Sat Nov 16 04:53:39 UTC 2024 with GitHub Actions
# Refresh Price / Volume(1h) figures, plus a re-ranking of a few coins
# (rows 43-46 swap which coin occupies which rank).
#
# Note: several "Price" values look numeric (e.g. 622.50, 0.999, 1.00) but
# must stay text cells (matching the source data, which keeps trailing
# zeros / fixed decimal places that a real number would lose). Prefixing
# those with a leading apostrophe forces Excel to store them as text
# instead of silently re-parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.280.84"
$ws.Range("E2").Value = "  +4.73%  "
$ws.Range("D3").Value = "3.131.77"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'220.47"
$ws.Range("E5").Value = "  +7.59%  "
$ws.Range("D6").Value = "'622.50"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +5.76%  "
$ws.Range("D8").Value = "'0.942"
$ws.Range("E8").Value = "  +17.84%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "3.130.07"
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("D11").Value = "'0.730"
$ws.Range("E11").Value = "  +26.27%  "
$ws.Range("E12").Value = "  +7.60%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  +10.58%  "
$ws.Range("D14").Value = "'34.62"
$ws.Range("E14").Value = "  +12.91%  "
$ws.Range("D15").Value = "'5.42"
$ws.Range("E15").Value = "  +4.96%  "
$ws.Range("D16").Value = "91.051.79"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").Value = "3.712.58"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("D18").Value = "3.130.50"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").Value = "'3.79"
$ws.Range("E19").Value = "  +21.81%  "
$ws.Range("E20").Value = "  +16.82%  "
$ws.Range("D21").Value = "'14.18"
$ws.Range("E21").Value = "  +11.21%  "
$ws.Range("D22").Value = "'432.66"
$ws.Range("E22").Value = "  +5.29%  "
$ws.Range("E23").Value = "  +10.71%  "
$ws.Range("E24").Value = "  +9.51%  "
$ws.Range("D25").Value = "'6.01"
$ws.Range("E25").Value = "  +15.52%  "
$ws.Range("D26").Value = "'12.38"
$ws.Range("E26").Value = "  +10.55%  "
$ws.Range("D27").Value = "'84.06"
$ws.Range("E27").Value = "  +5.90%  "
$ws.Range("D28").Value = "3.306.75"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +9.46%  "
$ws.Range("D31").Value = "'8.95"
$ws.Range("E31").Value = "  +13.91%  "
$ws.Range("D32").Value = "'529.39"
$ws.Range("E32").Value = "  +7.72%  "
$ws.Range("D33").Value = "'3.87"
$ws.Range("E33").Value = "  +14.54%  "
$ws.Range("D34").Value = "'0.884"
$ws.Range("E34").Value = "  -18.77%  "
$ws.Range("D35").Value = "'7.25"
$ws.Range("E35").Value = "  +13.58%  "
$ws.Range("E36").Value = "  +10.88%  "
$ws.Range("E37").Value = "  +9.71%  "
$ws.Range("E38").Value = "  +8.44%  "
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'0.151"
$ws.Range("E42").Value = "  +14.81%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.94"
$ws.Range("E44").Value = "  +10.78%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0770"
$ws.Range("E45").Value = "  +19.01%  "
$ws.Range("B46").Value = "PolygonEcosystemToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D46").Value = "'0.379"
$ws.Range("E46").Value = "  +8.34%  "
$ws.Range("D47").Value = "'144.44"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'44.19"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("E49").Value = "  +14.60%  "
$ws.Range("E50").Value = "  +29.75%  "
$ws.Range("D51").Value = "'168.00"
$ws.Range("E51").Value = "  +10.89%  "
